$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 796.8
$ws.Range("I20").Value = 796.8
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 796.8
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -566.8

$ws.Range("H35").Value = 796.8
$ws.Range("I35").Value = 796.8
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 796.8
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -417.8

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 42819336
$ws.Range("I74").Value = 80129570
$ws.Range("K74").Value = 80129570
$ws.Range("M74").Value = -80128696

$ws.Range("H77").Value = 42819336
$ws.Range("I77").Value = 80129570
$ws.Range("K77").Value = 400647850
$ws.Range("M77").Value = -400643482

$ws.Range("H122").Value = 2728.75
$ws.Range("I122").Value = 1638.3334
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 4915.0002
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -2465.0002
$ws.Range("N122").Value = -22900

$ws.Range("H132").Value = 41119910
$ws.Range("I132").Value = 40136124
$ws.Range("J132").Value = 42595590
$ws.Range("K132").Value = 120408372
$ws.Range("L132").Value = 127786770
$ws.Range("M132").Value = -120405842
$ws.Range("N132").Value = -127791830

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1955.56
$ws.Range("I86").Value = 1964.3469
$ws.Range("J86").Value = 1525
$ws.Range("K86").Value = 1964.3469
$ws.Range("L86").Value = 1525
$ws.Range("M86").Value = -841.3469
$ws.Range("N86").Value = -3771

$ws.Range("H89").Value = 1955.56
$ws.Range("I89").Value = 1964.3469
$ws.Range("J89").Value = 1525
$ws.Range("K89").Value = 9821.7345
$ws.Range("L89").Value = 7625
$ws.Range("M89").Value = -4205.7345
$ws.Range("N89").Value = -18857

$ws.Range("H134").Value = 35295132
$ws.Range("I134").Value = 50000900
$ws.Range("K134").Value = 150002700
$ws.Range("M134").Value = -150000165

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 55004.5
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 55004.5
$ws.Range("K15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = 55004.5
$ws.Range("N15").Value = -55344.5

$ws.Range("H22").Value = 307.66666
$ws.Range("I22").Value = 259.58334
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 259.58334
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 90.41665999999998
$ws.Range("N22").Value = -1200

$ws.Range("H36").Value = 19560.445
$ws.Range("I36").Value = 3200
$ws.Range("J36").Value = 32648.8
$ws.Range("K36").Value = 3200
$ws.Range("L36").Value = 32648.8
$ws.Range("M36").Value = -2812
$ws.Range("N36").Value = -33424.8

$ws.Range("H40").Value = 19560.445
$ws.Range("I40").Value = 3200
$ws.Range("J40").Value = 32648.8
$ws.Range("K40").Value = 3200
$ws.Range("L40").Value = 32648.8
$ws.Range("M40").Value = -3040
$ws.Range("N40").Value = -32968.8

$ws.Range("H122").Value = 13944
$ws.Range("I122").Value = 18364.8
$ws.Range("J122").Value = 2892
$ws.Range("K122").Value = 55094.39999999999
$ws.Range("L122").Value = 8676
$ws.Range("M122").Value = -52644.39999999999
$ws.Range("N122").Value = -13576

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 649
$ws.Range("I63").Value = 649
$ws.Range("K63").Value = 1947
$ws.Range("M63").Value = -1198

$ws.Range("H64").Value = 4081.4546
$ws.Range("I64").Value = 498
$ws.Range("J64").Value = 4877.778
$ws.Range("K64").Value = 1494
$ws.Range("L64").Value = 14633.334
$ws.Range("M64").Value = -1224
$ws.Range("N64").Value = -15173.334

$ws.Range("H66").Value = 649
$ws.Range("I66").Value = 649
$ws.Range("K66").Value = 5841
$ws.Range("M66").Value = -2097

$ws.Range("H67").Value = 4081.4546
$ws.Range("I67").Value = 498
$ws.Range("J67").Value = 4877.778
$ws.Range("K67").Value = 1494
$ws.Range("L67").Value = 14633.334
$ws.Range("M67").Value = -558
$ws.Range("N67").Value = -16505.334

$ws.Range("H68").Value = 1052.3334
$ws.Range("I68").Value = 676.1177
$ws.Range("J68").Value = 1280.75
$ws.Range("K68").Value = 2028.3531
$ws.Range("L68").Value = 3842.25
$ws.Range("M68").Value = -1217.3531
$ws.Range("N68").Value = -5464.25

$ws.Range("H70").Value = 2332.6667
$ws.Range("I70").Value = 998
$ws.Range("K70").Value = 2994
$ws.Range("M70").Value = -2679

$ws.Range("H71").Value = 1052.3334
$ws.Range("I71").Value = 676.1177
$ws.Range("J71").Value = 1280.75
$ws.Range("K71").Value = 6085.0593
$ws.Range("L71").Value = 11526.75
$ws.Range("M71").Value = -2029.0593
$ws.Range("N71").Value = -19638.75

$ws.Range("H73").Value = 2332.6667
$ws.Range("I73").Value = 998
$ws.Range("K73").Value = 2994
$ws.Range("M73").Value = -1902

$ws.Range("H75").Value = 1300
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 1420
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 4260
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -6256

$ws.Range("H78").Value = 1300
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 1420
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 12780
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -22764

$ws.Range("H112").Value = 4049.5
$ws.Range("I112").Value = 1526
$ws.Range("J112").Value = 5311.25
$ws.Range("K112").Value = 4578
$ws.Range("L112").Value = 15933.75
$ws.Range("M112").Value = -3470
$ws.Range("N112").Value = -18149.75

$ws.Range("H132").Value = 1699.7222
$ws.Range("J132").Value = 1678.9286
$ws.Range("L132").Value = 15110.3574
$ws.Range("N132").Value = -20170.3574

$ws.Range("H139").Value = 73037.57000000001
$ws.Range("I139").Value = 84371.664
$ws.Range("J139").Value = 5033
$ws.Range("K139").Value = 253114.992
$ws.Range("L139").Value = 15099
$ws.Range("M139").Value = -247974.992
$ws.Range("N139").Value = -25379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20165292
$ws.Range("I132").Value = 19167836
$ws.Range("J132").Value = 21253426
$ws.Range("K132").Value = 57503508
$ws.Range("L132").Value = 63760278
$ws.Range("M132").Value = -57500978
$ws.Range("N132").Value = -63765338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2754223
$ws.Range("I132").Value = 4280268.5
$ws.Range("J132").Value = 7340.8
$ws.Range("K132").Value = 12840805.5
$ws.Range("L132").Value = 22022.4
$ws.Range("M132").Value = -12838275.5
$ws.Range("N132").Value = -27082.4
